# DETC-104 - Updating section 3 - design 2 implementation details
#
# The "tRefContactType" helper table on Sheet2 (columns K:N, used to build a
# CREATE TABLE statement) had its "ContactTypeId" column/row removed - the
# table no longer declares that field, so every row below it moves up by
# one and the trailing "))" marker row is tidied up.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet2: tRefContactType (K:N) - drop the ContactTypeId row, shifting
# ContactTypeCode / ContactTypeDesc / IsActive / CreatedBy / CreatedDate /
# UpdatedBy / UpdatedDate up by one row. Column K (the literal separators
# "(" / ",") stays put - only L (field name), M (type) and N (the
# CONCATENATE formula) move.
# ---------------------------------------------------------------------

$ws2.Range("L3").Value = "ContactTypeCode"
$ws2.Range("M3").Value = " varchar(20)"
$ws2.Range("N3").Formula = "=CONCATENATE(K4,L3,M3)"
$ws2.Range("L3").ClearFormats()

$ws2.Range("L4").Value = "ContactTypeDesc"
$ws2.Range("M4").Value = " varchar(100)"
$ws2.Range("N4").Formula = "=CONCATENATE(K5,L4,M4)"

$ws2.Range("L5").Value = "IsActive"
$ws2.Range("M5").Value = " varchar(1)"
$ws2.Range("N5").Formula = "=CONCATENATE(K6,L5,M5)"

$ws2.Range("L6").Value = "CreatedBy"
$ws2.Range("M6").Value = " varchar(20)"
$ws2.Range("N6").Formula = "=CONCATENATE(K7,L6,M6)"

$ws2.Range("L7").Value = "CreatedDate"
$ws2.Range("M7").Value = " timestamp"
$ws2.Range("N7").Formula = "=CONCATENATE(K8,L7,M7)"

$ws2.Range("L8").Value = "UpdatedBy"
$ws2.Range("M8").Value = " varchar(20)"
$ws2.Range("N8").Formula = "=CONCATENATE(K9,L8,M8)"

$ws2.Range("L9").Value = "UpdatedDate"
$ws2.Range("M9").Value = " timestamp"
$ws2.Range("N9").Formula = "=CONCATENATE(K10,L9,M9)"

# Row 10 used to hold ",UpdatedDate timestamp"; the data rows now stop at
# row 9, so row 10 becomes just the closing ");" and no longer needs a
# K/L entry.
$ws2.Range("K10").ClearContents()
$ws2.Range("L10").ClearContents()
$ws2.Range("M10").Value = ");"
$ws2.Range("N10").Formula = "=CONCATENATE(K11,L10,M10)"

# Row 11's old closing ");" cells are gone entirely (the table now ends
# one row earlier).
$ws2.Range("M11").ClearContents()
$ws2.Range("N11").ClearContents()

# Row 13 (start of the next table, tRefWeightScale) becomes the new
# anchor/master of the shared CONCATENATE formula that used to live on N3.
$ws2.Range("N13").Formula = "=CONCATENATE(K13,L13,M13)"

# ---------------------------------------------------------------------
# Sheet1 / Sheet2 view selection, matching the saved cursor position.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("F24:G25").Select()

$ws2.Activate()
$ws2.Range("D40").Select()
